$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$durFormula = '=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])'

# Row 30
$ws.Range("E30").Value = 44267
$ws.Range("E30").NumberFormat = "m/d/yy"
$ws.Range("F30").Value = 0.57638888888888895
$ws.Range("F30").NumberFormat = "h:mm"
$ws.Range("G30").Value = 0.57986111111111105
$ws.Range("G30").NumberFormat = "h:mm"
$ws.Range("H30").Formula = $durFormula
$ws.Range("I30").Value = "Développement"
$ws.Range("J30").Value = "mise aux propre"
$ws.Range("K30").Value = "CPNV"
$ws.Range("L30").Value = "changement de commentaire, et mise au propre de certaine fonction"

# Row 31
$ws.Range("E31").Value = 44267
$ws.Range("E31").NumberFormat = "m/d/yy"
$ws.Range("F31").Value = 0.58333333333333337
$ws.Range("F31").NumberFormat = "h:mm"
$ws.Range("G31").Value = 0.58680555555555558
$ws.Range("G31").NumberFormat = "h:mm"
$ws.Range("H31").Formula = $durFormula
$ws.Range("I31").Value = "Développement"
$ws.Range("J31").Value = "Résolution de bug"
$ws.Range("K31").Value = "CPNV"
$ws.Range("L31").Value = "obligé l'utilisateur à entré les valeur correctes"

$ws.Range("L32").Select()
